$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend rows 3-5 with the same formatting as row 2 so newly written cells
# pick up the existing body-row style (font/alignment) instead of the
# workbook default style.
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A3:I5").PasteSpecial(-4122) | Out-Null

# --- roomId column (A): now a VALUE() formula shared down A2:A5 ---
$ws.Range("A2:A5").Formula = "=VALUE(55)"
$ws.Range("A2:A5").NumberFormat = "#,##0.00"

# --- roomName column (B): quote-prefixed numeric-looking text for the
#     first three rows, plain text for the last ---
$ws.Range("B2").Value = "'204"
$ws.Range("B3").Value = "'205"
$ws.Range("B4").Value = "'206"
$ws.Range("B5").Value = "not number"

# --- type column (C) ---
$ws.Range("C2").Value = "Single"
$ws.Range("C3").Value = "Double"
$ws.Range("C4").Value = "Suite"
$ws.Range("C5").Value = "Single"

# --- accessible column (D) ---
$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $false
$ws.Range("D4").Value = $false
$ws.Range("D5").Value = $true

# --- image column (E) ---
$ws.Range("E2").Value = "hiii"
$ws.Range("E3").Value = "hiii"
$ws.Range("E4").Value = "hiii"
$ws.Range("E5").Value = "hiii"

# --- description column (F) ---
$ws.Range("F2").Value = "Test desc1"
$ws.Range("F3").Value = "Test desc"
$ws.Range("F4").Value = "Test desc"
$ws.Range("F5").Value = "Test desc"

# --- features column (G) ---
$ws.Range("G2").Value = "WiFi|TV"
$ws.Range("G3").Value = "WiFi|TV"
$ws.Range("G4").Value = "WiFi|TV"
$ws.Range("G5").Value = "WiFi|TV"

# --- roomPrice column (H): VALUE() formula shared down H2:H5 ---
$ws.Range("H2:H5").Formula = "=VALUE(250)"
$ws.Range("H2:H5").NumberFormat = "#,##0.00"

# --- expected column (I) ---
$ws.Range("I2").Value = "pass"
$ws.Range("I3").Value = "pass"
$ws.Range("I4").Value = "pass"
$ws.Range("I5").Value = "pass"
